$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 4725.125
$ws.Range("I70").Value = 4883.5
$ws.Range("K70").Value = 14650.5
$ws.Range("M70").Value = -14380.5
# Row 73
$ws.Range("H73").Value = 4725.125
$ws.Range("I73").Value = 4883.5
$ws.Range("K73").Value = 14650.5
$ws.Range("M73").Value = -13714.5
# Row 137
$ws.Range("H137").Value = 2765.45
$ws.Range("I137").Value = 1115.2646
$ws.Range("J137").Value = 3985.152
$ws.Range("K137").Value = 3345.7938
$ws.Range("L137").Value = 11955.456
$ws.Range("M137").Value = -795.7937999999999
$ws.Range("N137").Value = -17055.456

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 34
$ws.Range("H34").Value = 900
$ws.Range("I34").Value = 900
$ws.Range("K34").Value = 900
$ws.Range("M34").Value = -629
# Row 88
$ws.Range("H88").Value = 2194
$ws.Range("I88").Value = 1504.6666
$ws.Range("J88").Value = 2607.6
$ws.Range("K88").Value = 1504.6666
$ws.Range("L88").Value = 2607.6
$ws.Range("M88").Value = -1098.6666
$ws.Range("N88").Value = -3419.6
# Row 91
$ws.Range("H91").Value = 2194
$ws.Range("I91").Value = 1504.6666
$ws.Range("J91").Value = 2607.6
$ws.Range("K91").Value = 1504.6666
$ws.Range("L91").Value = 2607.6
$ws.Range("M91").Value = -100.6666
$ws.Range("N91").Value = -5415.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 36
$ws.Range("H36").Value = 5418.222
$ws.Range("I36").Value = 680.5714
$ws.Range("J36").Value = 22000
$ws.Range("K36").Value = 680.5714
$ws.Range("L36").Value = 22000
$ws.Range("M36").Value = -146.5714
$ws.Range("N36").Value = -23068
# Row 123
$ws.Range("H123").Value = 40508.11
$ws.Range("J123").Value = 40508.11
$ws.Range("L123").Value = 40508.11
$ws.Range("N123").Value = -50308.11

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 846.5
$ws.Range("I4").Value = 224.44444
$ws.Range("J4").Value = 1157.5278
$ws.Range("K4").Value = 673.33332
$ws.Range("L4").Value = 3472.5834
$ws.Range("M4").Value = -561.33332
$ws.Range("N4").Value = -3696.5834
# Row 17
$ws.Range("H17").Value = 8713.5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 8713.5
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 26140.5
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -26478.5
# Row 22
$ws.Range("H22").Value = 67067332
$ws.Range("I22").Value = 666770.3
$ws.Range("J22").Value = 166668180
$ws.Range("K22").Value = 2000310.9
$ws.Range("L22").Value = 500004540
$ws.Range("M22").Value = -2000141.9
$ws.Range("N22").Value = -500004878
# Row 25
$ws.Range("H25").Value = 290
$ws.Range("I25").Value = 290
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 870
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -701
$ws.Range("N25").ClearContents()
# Row 27
$ws.Range("H27").Value = 67067332
$ws.Range("I27").Value = 666770.3
$ws.Range("J27").Value = 166668180
$ws.Range("K27").Value = 2000310.9
$ws.Range("L27").Value = 500004540
$ws.Range("M27").Value = -2000208.9
$ws.Range("N27").Value = -500004744
# Row 30
$ws.Range("H30").Value = 290
$ws.Range("I30").Value = 290
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 870
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -768
$ws.Range("N30").ClearContents()
# Row 32
$ws.Range("H32").Value = 414285.72
$ws.Range("I32").Value = 414285.72
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1242857.16
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1242574.16
$ws.Range("N32").ClearContents()
# Row 37
$ws.Range("H37").Value = 40000
$ws.Range("J37").Value = 40000
$ws.Range("L37").Value = 120000
$ws.Range("N37").Value = -120224
# Row 49
$ws.Range("H49").Value = 3000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
# Row 64
$ws.Range("H64").Value = 111113120
$ws.Range("I64").Value = 1012
$ws.Range("J64").Value = 166669170
$ws.Range("K64").Value = 3036
$ws.Range("L64").Value = 500007510
$ws.Range("M64").Value = -2766
$ws.Range("N64").Value = -500008050
# Row 67
$ws.Range("H67").Value = 111113120
$ws.Range("I67").Value = 1012
$ws.Range("J67").Value = 166669170
$ws.Range("K67").Value = 3036
$ws.Range("L67").Value = 500007510
$ws.Range("M67").Value = -2100
$ws.Range("N67").Value = -500009382
# Row 70
$ws.Range("H70").Value = 3680.75
$ws.Range("I70").Value = 3206.5715
$ws.Range("J70").Value = 7000
$ws.Range("K70").Value = 9619.7145
$ws.Range("L70").Value = 21000
$ws.Range("M70").Value = -9304.7145
$ws.Range("N70").Value = -21630
# Row 73
$ws.Range("H73").Value = 3680.75
$ws.Range("I73").Value = 3206.5715
$ws.Range("J73").Value = 7000
$ws.Range("K73").Value = 9619.7145
$ws.Range("L73").Value = 21000
$ws.Range("M73").Value = -8527.7145
$ws.Range("N73").Value = -23184
# Row 94
$ws.Range("H94").Value = 1875
$ws.Range("I94").Value = 625
$ws.Range("J94").Value = 2500
$ws.Range("K94").Value = 1875
$ws.Range("L94").Value = 7500
$ws.Range("M94").Value = -1199
$ws.Range("N94").Value = -8852
# Row 95
$ws.Range("H95").Value = 3166.6667
$ws.Range("I95").Value = 1500
$ws.Range("J95").Value = 4000
$ws.Range("K95").Value = 4500
$ws.Range("L95").Value = 12000
$ws.Range("M95").Value = -2441
$ws.Range("N95").Value = -16118
# Row 97
$ws.Range("H97").Value = 910.75
$ws.Range("I97").Value = 547.6667
$ws.Range("K97").Value = 1643.0001
$ws.Range("M97").Value = -1147.0001
# Row 100
$ws.Range("H100").Value = 7998
$ws.Range("J100").Value = 7998
$ws.Range("L100").Value = 23994
$ws.Range("N100").Value = -25616
# Row 103
$ws.Range("H103").Value = 688.6667
$ws.Range("I103").Value = 688.6667
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 2066.0001
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -1187.0001
$ws.Range("N103").ClearContents()
# Row 106
$ws.Range("H106").Value = 145828.42
$ws.Range("J106").Value = 145828.42
$ws.Range("L106").Value = 437485.26
$ws.Range("N106").Value = -439377.26

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 4375.3403
$ws.Range("I132").Value = 4568.1353
$ws.Range("J132").Value = 3662
$ws.Range("K132").Value = 13704.4059
$ws.Range("L132").Value = 10986
$ws.Range("M132").Value = -11174.4059
$ws.Range("N132").Value = -16046
